$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 4 values / formulas (B4, C4, D4, E4 were empty before)
$ws.Range("B4").Value = 32
$ws.Range("C4").Formula = "=(8.67+8.21)/2"
$ws.Range("D4").Formula = "=1"
$ws.Range("E4").Value = 7.5

# Move the active selection to G5 (was I7)
$ws.Range("G5").Select()
